$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 157, pushing existing rows 157-168 down to 159-170.
$ws.Range("A157:A158").EntireRow.Insert()

# Populate new row 157
$ws.Cells.Item(157,1).Value = 7
$ws.Cells.Item(157,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(157,3).Value = "Ñuble"
$ws.Cells.Item(157,4).Value = 44610
$ws.Cells.Item(157,5).Value = 16
$ws.Cells.Item(157,6).Value = 100112024
$ws.Cells.Item(157,7).Value = "Choclo"
$ws.Cells.Item(157,8).Value = "Choclero"
$ws.Cells.Item(157,9).Value = "Primera"
$ws.Cells.Item(157,10).Value = 12000
$ws.Cells.Item(157,11).Value = 150
$ws.Cells.Item(157,12).Value = 180
$ws.Cells.Item(157,13).Value = 165
$ws.Cells.Item(157,14).Value = "`$/unidad"
$ws.Cells.Item(157,15).Value = "Región del Maule"
$ws.Cells.Item(157,16).Value = 165
$ws.Cells.Item(157,17).Value = 1
$ws.Cells.Item(157,18).Value = "Hortaliza"

# Populate new row 158
$ws.Cells.Item(158,1).Value = 7
$ws.Cells.Item(158,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(158,3).Value = "Ñuble"
$ws.Cells.Item(158,4).Value = 44610
$ws.Cells.Item(158,5).Value = 16
$ws.Cells.Item(158,6).Value = 100112024
$ws.Cells.Item(158,7).Value = "Choclo"
$ws.Cells.Item(158,8).Value = "Choclero"
$ws.Cells.Item(158,9).Value = "Segunda"
$ws.Cells.Item(158,10).Value = 10000
$ws.Cells.Item(158,11).Value = 100
$ws.Cells.Item(158,12).Value = 120
$ws.Cells.Item(158,13).Value = 110
$ws.Cells.Item(158,14).Value = "`$/unidad"
$ws.Cells.Item(158,15).Value = "Región del Maule"
$ws.Cells.Item(158,16).Value = 110
$ws.Cells.Item(158,17).Value = 1
$ws.Cells.Item(158,18).Value = "Hortaliza"

Write-Output "done"
